$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of rows 2-4 (full records) by one position:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 3
# Only the cells whose value actually differs between the old and new
# row contents are written below (matching the supplied diff exactly).

# --- Row 2 (becomes former row 4's record) ---
$ws.Range("A2").Value = 104023398
$ws.Range("Q2").Value = 543928.1404005223
$ws.Range("R2").Value = 7094070.257953409
$ws.Range("Z2").Value = "10:59"
$ws.Range("AB2").Value = "10:59"
$ws.Range("AO2").Value = "Sälg"

# --- Row 3 (becomes former row 2's record) ---
$ws.Range("A3").Value = 104023416
$ws.Range("B3").Value = 78569
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6458
$ws.Range("F3").Value = "Lunglav"
$ws.Range("G3").Value = "Lobaria pulmonaria"
$ws.Range("H3").Value = "(L.) Hoffm."
$ws.Range("Q3").Value = 543743.4593651614
$ws.Range("R3").Value = 7094441.509715242
$ws.Range("Z3").Value = "08:25"
$ws.Range("AB3").Value = "08:25"

# --- Row 4 (becomes former row 3's record) ---
$ws.Range("A4").Value = 104023407
$ws.Range("B4").Value = 78458
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 6457
$ws.Range("F4").Value = "Dvärgtufs"
$ws.Range("G4").Value = "Scytinium teretiusculum"
$ws.Range("H4").Value = "(Wallr.) Otálora et al."
$ws.Range("Q4").Value = 543798.3199733114
$ws.Range("R4").Value = 7094387.369215799
$ws.Range("Z4").Value = "09:55"
$ws.Range("AB4").Value = "09:55"
$ws.Range("AO4").Value = "Asp"
